$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update inputs: Bar Length (B2) and Grip Length (B3)
$ws.Range("B2").Value = 0.35
$ws.Range("B3").Value = 0.18

# Update Avg Lift Force (M4)
$ws.Range("M4").Value = 15

# Update the active selection on the sheet
$ws.Range("G10").Select() | Out-Null
